$d = $word.ActiveDocument

# The document currently ends with an empty paragraph (sz=28) right before
# the sectPr. We append five new paragraphs after it:
#   1. "Date: 24/1/25"                              (bold, 15pt)
#   2. "Topics:"                                     (bold, 14pt)
#   3. "COMPLETED DBMS Task – 2."                     (list item, numId 6)
#   4. "Design Database for ecommerce website."       (list item, numId 6)
#   5. "Partitioning and CTEs."                       (list item, numId 6)

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaIndex)

# Create a new paragraph after the existing trailing empty paragraph, then
# collapse its range to the very start so InsertXML lands exactly there
# without consuming/merging the original empty paragraph.
$lastPara.Range.InsertParagraphAfter()

$newParaIndex = $d.Paragraphs.Count
$insertRange = $d.Paragraphs($newParaIndex).Range
$insertRange.Collapse(1)

$xml = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Date: </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="en-US"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="en-US"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="en-US"/></w:rPr><w:t>/1/25</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Topics:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>C</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">OMPLETED </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">DBMS Task – </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>2.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Design Database for ecommerce website.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Partitioning and CTEs.</w:t></w:r></w:p>'

$insertRange.InsertXML($xml)

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
